$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 183-186 down to 184-187 (values only; the per-column cell
# styles of the destination rows already match the source rows because
# whole formatting "blocks" move together), working from the bottom up
# so we never clobber data before it has been copied.

# old row 186 (more_ims) -> row 187
$ws.Range("A187").Value2 = $ws.Range("A186").Value2()
$ws.Range("B187").Value2 = $ws.Range("B186").Value2()
$ws.Range("C187").Value2 = $ws.Range("C186").Value2()
$ws.Range("D187").Value2 = $ws.Range("D186").Value2()
$ws.Range("E187").Value2 = $ws.Range("E186").Value2()
$ws.Range("F187").Value2 = $ws.Range("F186").Value2()
$ws.Range("G187").Value2 = $ws.Range("G186").Value2()
$ws.Range("H187").Value2 = $ws.Range("H186").Value2()
$ws.Range("I187").Value2 = $ws.Range("I186").Value2()

# old row 185 (less_ims) -> row 186
$ws.Range("A186").Value2 = $ws.Range("A185").Value2()
$ws.Range("B186").Value2 = $ws.Range("B185").Value2()
$ws.Range("C186").Value2 = $ws.Range("C185").Value2()
$ws.Range("D186").Value2 = $ws.Range("D185").Value2()
$ws.Range("E186").Value2 = $ws.Range("E185").Value2()
$ws.Range("F186").Value2 = $ws.Range("F185").Value2()
$ws.Range("G186").Value2 = $ws.Range("G185").Value2()
$ws.Range("H186").Value2 = $ws.Range("H185").Value2()
$ws.Range("I186").Value2 = $ws.Range("I185").Value2()

# old row 184 (empty_private) -> row 185
$ws.Range("A185").Value2 = $ws.Range("A184").Value2()
$ws.Range("B185").Value2 = $ws.Range("B184").Value2()
$ws.Range("C185").Value2 = $ws.Range("C184").Value2()
$ws.Range("D185").Value2 = $ws.Range("D184").Value2()
$ws.Range("E185").Value2 = $ws.Range("E184").Value2()
$ws.Range("F185").Value2 = $ws.Range("F184").Value2()
$ws.Range("G185").Value2 = $ws.Range("G184").Value2()
$ws.Range("H185").Value2 = $ws.Range("H184").Value2()
$ws.Range("I185").Value2 = $ws.Range("I184").Value2()

# old row 183 (empty_coop) -> row 184
$ws.Range("A184").Value2 = $ws.Range("A183").Value2()
$ws.Range("B184").Value2 = $ws.Range("B183").Value2()
$ws.Range("C184").Value2 = $ws.Range("C183").Value2()
$ws.Range("D184").Value2 = $ws.Range("D183").Value2()
$ws.Range("E184").Value2 = $ws.Range("E183").Value2()
$ws.Range("F184").Value2 = $ws.Range("F183").Value2()
$ws.Range("G184").Value2 = $ws.Range("G183").Value2()
$ws.Range("H184").Value2 = $ws.Range("H183").Value2()
$ws.Range("I184").Value2 = $ws.Range("I183").Value2()

# new row 183: car_trust parameter (keeps the existing formatting of this
# block of rows, only the values change)
$ws.Range("A183").Value2 = "housing model"
$ws.Range("B183").Value2 = "car_trust"
$ws.Range("C183").Value2 = "no parameter in previous model"
$ws.Range("D183").Value2 = 80
$ws.Range("E183").Value2 = 80
$ws.Range("F183").Value2 = 80
$ws.Range("G183").Value2 = "percent"
$ws.Range("H183").Value2 = "medium"
$ws.Range("I183").Value2 = "If the amount of people according to the project list exceeds the population according to capacity/reserves: trust the capacity/reserves number (parameter = 100%)? Or the project list (parameter = 0%)?"

# Row heights: row 183 is a new wrapped-text row (same height class as the
# other rows in this block), rows 186/187 keep the height of the
# "demography and housing model" rows they inherited from 185/186.
$ws.Rows.Item(183).RowHeight = 38.25
$ws.Rows.Item(186).RowHeight = 51
$ws.Rows.Item(187).RowHeight = 51

# Scroll the frozen pane down so row 167 is the first visible data row.
$ws.Range("A167").Select()
